# v0.19 edit script
# Renames the two "sala oscura" rooms to "habitación de reclusión A" / "B"
# in three places in the guide text.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "Ir a sala oscura 2 (...)" -> "Ir a habitación de reclusión A (...)"
# The original run "sala oscura 2" is immediately followed by a run that
# only contains a single trailing space; replacing the whole "sala oscura 2 "
# span collapses both into the new single run, which matches the target.
# ---------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("sala oscura 2 ", $false, $false, $false, $false, $false, $true, 1, $false, "habitación de reclusión A ", 2)
if (-not $found1) {
    Write-Output "WARNING: change 1 target not found"
}

# ---------------------------------------------------------------------
# Change 2: "...el correcto en la sala oscura 1 (...)" ->
#           "...el correcto en la habitación de reclusión B (...)"
# Here the original keeps "la" in its own run and " " in the run after it,
# and the target keeps that exact same two-run split ("la" / " habitación
# de reclusión B "). A plain Find/Replace across the run boundary would
# coalesce both spans into a single run, so we drop a temporary zero-width
# bookmark between them to force the writer to keep them as two runs, then
# remove the bookmark once the text has been updated.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("la sala oscura 1")
if (-not $found2) {
    Write-Output "WARNING: change 2 target not found"
} else {
    $cut2 = $rng2.Start + 2   # just after "la"
    $wall2 = $d.Range($cut2, $cut2)
    $d.Bookmarks.Add("zzWall2", $wall2) | Out-Null

    $spaceRng = $d.Range($d.Bookmarks("zzWall2").Range.Start, $d.Bookmarks("zzWall2").Range.Start + 1)
    $spaceRng.Text = " habitación de reclusión B "

    $d.Bookmarks("zzWall2").Delete()
}

# ---------------------------------------------------------------------
# Change 3: "...La sala oscura 1 se abre con la llave alfil." ->
#           "...La habitación de reclusión B se abre con la llave alfil."
# The target splits this sentence into three runs ("La ", "habitación de
# reclusión B ", "se abre con la llave alfil") that all share the same
# formatting, so again we use temporary bookmarks as walls to stop the
# writer from re-merging them, then delete the bookmarks.
# ---------------------------------------------------------------------
$rng3 = $d.Content
$found3 = $rng3.Find.Execute("La sala oscura 1 se abre con la llave alfil")
if (-not $found3) {
    Write-Output "WARNING: change 3 target not found"
} else {
    $base3 = $rng3.Start
    $cutA = $base3 + 3     # just after "La "
    $cutB = $base3 + 17    # just after "sala oscura 1 "

    $wallA = $d.Range($cutA, $cutA)
    $d.Bookmarks.Add("zzWallA", $wallA) | Out-Null
    $wallB = $d.Range($cutB, $cutB)
    $d.Bookmarks.Add("zzWallB", $wallB) | Out-Null

    $midRng = $d.Range($d.Bookmarks("zzWallA").Range.Start, $d.Bookmarks("zzWallB").Range.Start)
    $midRng.Text = "habitación de reclusión B "

    $d.Bookmarks("zzWallA").Delete()
    $d.Bookmarks("zzWallB").Delete()
}

Write-Output "done"
